$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows to append after the last existing row (183)
$newData = @(
    @("20-09-2021", 3.47, 3.18, 3.3),
    @("21-09-2021", 3.42, 3.13, 3.29),
    @("22-09-2021", 3.47, 3.15, 3.28),
    @("23-09-2021", 3.32, 3.18, 3.27),
    @("24-09-2021", 3.24, 3.19, 3.24)
)

$startRow = 184
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $values = $newData[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}
